# Update column C ("Förändrad") values for rows 2-36 from 45688 to 45690
# (i.e. from 2025-01-31 to 2025-02-02), keeping existing cell formatting.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

for ($row = 2; $row -le 36; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45688) {
        $cell.Value = 45690
    }
}
